$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 37.98277566666666
$ws.Range("H2").Value = 113.948327
$ws.Range("I2").Value = 0.697850645410475
$ws.Range("J2").Value = 0.6978506454104751
$ws.Range("M2").Value = 1.118034
$ws.Range("N2").Value = 3.354102
$ws.Range("O2").Value = 0.01817106018145251
$ws.Range("P2").Value = 0.01817106018145251
$ws.Range("Q2").Value = 42.46603460970599
$ws.Range("R2").Value = 382.194311487354
$ws.Range("S2").Value = 0.01268068607541922
$ws.Range("T2").Value = 0.01268068607541922
$ws.Range("G3").Value = 37.98277566666666
$ws.Range("H3").Value = 113.948327
$ws.Range("I3").Value = 0.697850645410475
$ws.Range("J3").Value = 0.6978506454104751
$ws.Range("O3").Value = 0.962887615892719
$ws.Range("P3").Value = 0.9628876158927191
$ws.Range("Q3").Value = 2250.282504897243
$ws.Range("R3").Value = 20252.54254407519
$ws.Range("S3").Value = 0.6719517442084875
$ws.Range("T3").Value = 0.6719517442084877
$ws.Range("G4").Value = 37.98277566666666
$ws.Range("H4").Value = 113.948327
$ws.Range("I4").Value = 0.697850645410475
$ws.Range("J4").Value = 0.6978506454104751
$ws.Range("M4").Value = 1.121724666666667
$ws.Range("N4").Value = 3.365174
$ws.Range("O4").Value = 0.01823104344324033
$ws.Range("P4").Value = 0.01823104344324033
$ws.Range("Q4").Value = 42.60621637376644
$ws.Range("R4").Value = 383.455947363898
$ws.Range("S4").Value = 0.01272254543337167
$ws.Range("T4").Value = 0.01272254543337167
$ws.Range("A5").Value = "ECs"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 37.98277566666666
$ws.Range("H5").Value = 113.948327
$ws.Range("I5").Value = 0.697850645410475
$ws.Range("J5").Value = 0.6978506454104751
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04370233333333334
$ws.Range("N5").Value = 0.131107
$ws.Range("O5").Value = 0.0007102804825880949
$ws.Range("P5").Value = 0.0007102804825880949
$ws.Range("Q5").Value = 1.659935923109889
$ws.Range("R5").Value = 14.939423307989
$ws.Range("S5").Value = 0.0004956696931965657
$ws.Range("T5").Value = 0.0004956696931965658
$ws.Range("D6").Value = "ECs"
$ws.Range("I6").Value = 0.1779541659542351
$ws.Range("J6").Value = 0.1779541659542352
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.118034
$ws.Range("N6").Value = 3.354102
$ws.Range("O6").Value = 0.01817106018145251
$ws.Range("P6").Value = 0.01817106018145251
$ws.Range("Q6").Value = 10.828975827498
$ws.Range("R6").Value = 97.46078244748202
$ws.Range("S6").Value = 0.003233615859094594
$ws.Range("T6").Value = 0.003233615859094595
$ws.Range("D7").Value = "FAPs"
$ws.Range("I7").Value = 0.1779541659542351
$ws.Range("J7").Value = 0.1779541659542352
$ws.Range("M7").Value = 59.24481466666666
$ws.Range("N7").Value = 177.734444
$ws.Range("O7").Value = 0.962887615892719
$ws.Range("P7").Value = 0.9628876158927191
$ws.Range("Q7").Value = 573.8292985096449
$ws.Range("R7").Value = 5164.463686586804
$ws.Range("S7").Value = 0.1713498625938507
$ws.Range("T7").Value = 0.1713498625938508
$ws.Range("A8").Value = "FAPs"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("G8").Value = 9.685730333333334
$ws.Range("H8").Value = 29.057191
$ws.Range("I8").Value = 0.1779541659542351
$ws.Range("J8").Value = 0.1779541659542352
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.121724666666667
$ws.Range("N8").Value = 3.365174
$ws.Range("O8").Value = 0.01823104344324033
$ws.Range("P8").Value = 0.01823104344324033
$ws.Range("Q8").Value = 10.86472262958156
$ws.Range("R8").Value = 97.78250366623401
$ws.Range("S8").Value = 0.00324429013041726
$ws.Range("T8").Value = 0.00324429013041726
$ws.Range("A9").Value = "FAPs"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("G9").Value = 9.685730333333334
$ws.Range("H9").Value = 29.057191
$ws.Range("I9").Value = 0.1779541659542351
$ws.Range("J9").Value = 0.1779541659542352
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.04370233333333334
$ws.Range("N9").Value = 0.131107
$ws.Range("O9").Value = 0.0007102804825880949
$ws.Range("P9").Value = 0.0007102804825880949
$ws.Range("Q9").Value = 0.4232890156041111
$ws.Range("R9").Value = 3.809601140437
$ws.Range("S9").Value = 0.0001263973708725361
$ws.Range("T9").Value = 0.0001263973708725361
$ws.Range("D10").Value = "ECs"
$ws.Range("G10").Value = 0.5676613333333332
$ws.Range("H10").Value = 1.702984
$ws.Range("I10").Value = 0.01042953867610283
$ws.Range("J10").Value = 0.01042953867610283
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 1.118034
$ws.Range("N10").Value = 3.354102
$ws.Range("O10").Value = 0.01817106018145251
$ws.Range("P10").Value = 0.01817106018145251
$ws.Range("Q10").Value = 0.6346646711519999
$ws.Range("R10").Value = 5.711982040368
$ws.Range("S10").Value = 0.000189515774948251
$ws.Range("T10").Value = 0.0001895157749482511
$ws.Range("A11").Value = "Inflammatory-Mac"
$ws.Range("D11").Value = "FAPs"
$ws.Range("G11").Value = 0.5676613333333332
$ws.Range("H11").Value = 1.702984
$ws.Range("I11").Value = 0.01042953867610283
$ws.Range("J11").Value = 0.01042953867610283
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 59.24481466666666
$ws.Range("N11").Value = 177.734444
$ws.Range("O11").Value = 0.962887615892719
$ws.Range("P11").Value = 0.9628876158927191
$ws.Range("Q11").Value = 33.63099048676622
$ws.Range("R11").Value = 302.6789143808959
$ws.Range("S11").Value = 0.01004247363069356
$ws.Range("T11").Value = 0.01004247363069356
$ws.Range("A12").Value = "Inflammatory-Mac"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("G12").Value = 0.5676613333333332
$ws.Range("H12").Value = 1.702984
$ws.Range("I12").Value = 0.01042953867610283
$ws.Range("J12").Value = 0.01042953867610283
$ws.Range("M12").Value = 1.121724666666667
$ws.Range("N12").Value = 3.365174
$ws.Range("O12").Value = 0.01823104344324033
$ws.Range("P12").Value = 0.01823104344324033
$ws.Range("Q12").Value = 0.6367597199128888
$ws.Range("R12").Value = 5.730837479216
$ws.Range("S12").Value = 0.0001901413726969859
$ws.Range("T12").Value = 0.000190141372696986
$ws.Range("A13").Value = "Inflammatory-Mac"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("G13").Value = 0.5676613333333332
$ws.Range("H13").Value = 1.702984
$ws.Range("I13").Value = 0.01042953867610283
$ws.Range("J13").Value = 0.01042953867610283
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.04370233333333334
$ws.Range("N13").Value = 0.131107
$ws.Range("O13").Value = 0.0007102804825880949
$ws.Range("P13").Value = 0.0007102804825880949
$ws.Range("Q13").Value = 0.02480812480977778
$ws.Range("R13").Value = 0.223273123288
$ws.Range("S13").Value = (7.407897764033519 / 1000000)
$ws.Range("T13").Value = (7.40789776403352 / 1000000)
$ws.Range("A14").Value = "MuSCs"
$ws.Range("G14").Value = 5.823095333333334
$ws.Range("H14").Value = 17.469286
$ws.Range("I14").Value = 0.1069866739681064
$ws.Range("J14").Value = 0.1069866739681064
$ws.Range("M14").Value = 1.118034
$ws.Range("N14").Value = 3.354102
$ws.Range("O14").Value = 0.01817106018145251
$ws.Range("P14").Value = 0.01817106018145251
$ws.Range("Q14").Value = 6.510418567908
$ws.Range("R14").Value = 58.59376711117201
$ws.Range("S14").Value = 0.0019440612912879
$ws.Range("T14").Value = 0.001944061291287901
$ws.Range("A15").Value = "MuSCs"
$ws.Range("G15").Value = 5.823095333333334
$ws.Range("H15").Value = 17.469286
$ws.Range("I15").Value = 0.1069866739681064
$ws.Range("J15").Value = 0.1069866739681064
$ws.Range("O15").Value = 0.962887615892719
$ws.Range("P15").Value = 0.9628876158927191
$ws.Range("Q15").Value = 344.9882038096649
$ws.Range("R15").Value = 3104.893834286984
$ws.Range("S15").Value = 0.1030161434294416
$ws.Range("T15").Value = 0.1030161434294416
$ws.Range("A16").Value = "MuSCs"
$ws.Range("G16").Value = 5.823095333333334
$ws.Range("H16").Value = 17.469286
$ws.Range("I16").Value = 0.1069866739681064
$ws.Range("J16").Value = 0.1069866739681064
$ws.Range("M16").Value = 1.121724666666667
$ws.Range("N16").Value = 3.365174
$ws.Range("O16").Value = 0.01823104344324033
$ws.Range("P16").Value = 0.01823104344324033
$ws.Range("Q16").Value = 6.531909671751556
$ws.Range("R16").Value = 58.787187045764
$ws.Range("S16").Value = 0.001950478700960338
$ws.Range("T16").Value = 0.001950478700960338
$ws.Range("A17").Value = "MuSCs"
$ws.Range("B17").Value = "Efnb2"
$ws.Range("C17").Value = "Epha3"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 5.823095333333334
$ws.Range("H17").Value = 17.469286
$ws.Range("I17").Value = 0.1069866739681064
$ws.Range("J17").Value = 0.1069866739681064
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.04370233333333334
$ws.Range("N17").Value = 0.131107
$ws.Range("O17").Value = 0.0007102804825880949
$ws.Range("P17").Value = 0.0007102804825880949
$ws.Range("Q17").Value = 0.2544828532891111
$ws.Range("R17").Value = 2.290345679602
$ws.Range("S17").Value = (7.59905464165618 / 100000)
$ws.Range("T17").Value = (7.59905464165618 / 100000)
$ws.Range("A18").Value = "Resolving-Mac"
$ws.Range("B18").Value = "Efnb2"
$ws.Range("C18").Value = "Epha3"
$ws.Range("D18").Value = "ECs"
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.3689676666666666
$ws.Range("H18").Value = 1.106903
$ws.Range("I18").Value = 0.006778975991080511
$ws.Range("J18").Value = 0.006778975991080512
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 0.6666666666666666
$ws.Range("M18").Value = 1.118034
$ws.Range("N18").Value = 3.354102
$ws.Range("O18").Value = 0.01817106018145251
$ws.Range("P18").Value = 0.01817106018145251
$ws.Range("Q18").Value = 0.4125183962339999
$ws.Range("R18").Value = 3.712665566106
$ws.Range("S18").Value = 0.0001231811807025456
$ws.Range("T18").Value = 0.0001231811807025457
$ws.Range("A19").Value = "Resolving-Mac"
$ws.Range("B19").Value = "Efnb2"
$ws.Range("C19").Value = "Epha3"
$ws.Range("D19").Value = "FAPs"
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.3689676666666666
$ws.Range("H19").Value = 1.106903
$ws.Range("I19").Value = 0.006778975991080511
$ws.Range("J19").Value = 0.006778975991080512
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 59.24481466666666
$ws.Range("N19").Value = 177.734444
$ws.Range("O19").Value = 0.962887615892719
$ws.Range("P19").Value = 0.9628876158927191
$ws.Range("Q19").Value = 21.85942102965911
$ws.Range("R19").Value = 196.734789266932
$ws.Range("S19").Value = 0.006527392030245494
$ws.Range("T19").Value = 0.006527392030245496
$ws.Range("A20").Value = "Resolving-Mac"
$ws.Range("B20").Value = "Efnb2"
$ws.Range("C20").Value = "Epha3"
$ws.Range("D20").Value = "MuSCs"
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.3689676666666666
$ws.Range("H20").Value = 1.106903
$ws.Range("I20").Value = 0.006778975991080511
$ws.Range("J20").Value = 0.006778975991080512
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 1.121724666666667
$ws.Range("N20").Value = 3.365174
$ws.Range("O20").Value = 0.01823104344324033
$ws.Range("P20").Value = 0.01823104344324033
$ws.Range("Q20").Value = 0.4138801329024444
$ws.Range("R20").Value = 3.724921196122
$ws.Range("S20").Value = 0.000123587805794072
$ws.Range("T20").Value = 0.000123587805794072
$ws.Range("A21").Value = "Resolving-Mac"
$ws.Range("B21").Value = "Efnb2"
$ws.Range("C21").Value = "Epha3"
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.3689676666666666
$ws.Range("H21").Value = 1.106903
$ws.Range("I21").Value = 0.006778975991080511
$ws.Range("J21").Value = 0.006778975991080512
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 0.3333333333333333
$ws.Range("M21").Value = 0.04370233333333334
$ws.Range("N21").Value = 0.131107
$ws.Range("O21").Value = 0.0007102804825880949
$ws.Range("P21").Value = 0.0007102804825880949
$ws.Range("Q21").Value = 0.01612474795788889
$ws.Range("R21").Value = 0.145122731621
$ws.Range("S21").Value = (4.814974338397774 / 1000000)
$ws.Range("T21").Value = (4.814974338397775 / 1000000)